{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" footer paragraph that, together with the\n// \"\u00a9 2020 ...\" paragraph right after it and the blank paragraph right\n// before it, should be removed (the bibliography entry should then be\n// followed directly by the existing blank / page-break paragraphs).\nconst targetText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst idx = items.findIndex((p) => p.text === targetText);\n\nif (idx !== -1) {\n  // Delete the \"\u00a9 2020 ...\" paragraph that follows the footer line.\n  if (idx + 1 < items.length) {\n    items[idx + 1].delete();\n  }\n  // Delete the \"Ver no Jupiter ...\" footer line itself.\n  items[idx].delete();\n  // Delete the blank paragraph that separated it from the bibliography text.\n  if (idx - 1 >= 0) {\n    items[idx - 1].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer paragraphs\n# and the blank paragraph that separated them from the bibliography text,\n# leaving the bibliography entry directly followed by the existing blank\n# paragraph and the page-break paragraph.\n$target = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {\n        # Delete this paragraph, the \"\u00a9 2020 ...\" paragraph right after it,\n        # and the blank paragraph right before it (3 paragraphs total).\n        $d.Paragraphs.Item($i + 1).Range.Delete()\n        $d.Paragraphs.Item($i).Range.Delete()\n        $d.Paragraphs.Item($i - 1).Range.Delete()\n        break\n    }\n}\n"}
